$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Contest 48: RCB vs PBKS (row 60) - player scores
$ws.Range("E60").Value = 40
$ws.Range("H60").Value = 30
$ws.Range("K60").Value = 80
$ws.Range("N60").Value = 50
$ws.Range("Q60").Value = 100
$ws.Range("T60").Value = 60
$ws.Range("W60").Value = 20
$ws.Range("Z60").Value = 0
$ws.Range("AC60").Value = 70

# Contest 49: KKR vs SRH (row 61) - player scores
$ws.Range("E61").Value = 40
$ws.Range("H61").Value = 20
$ws.Range("K61").Value = 80
$ws.Range("N61").Value = 0
$ws.Range("Q61").Value = 30
$ws.Range("T61").Value = 60
$ws.Range("W61").Value = 100
$ws.Range("Z61").Value = 70
$ws.Range("AC61").Value = 50
